$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently ends at row 64 with two Cilantro "Primera" price records
# (dated 2021-03-11 / value 44266, and 2022-02-17 / value 44609).
# A new week of data (2022-08-03 / value 44776) is being added as the new
# top two rows (63 & 64: "Primera" + "Segunda"), pushing the two existing
# rows down to become rows 65 & 66.

# 1) Duplicate the existing row 63 down to become the new row 65
#    (preserves values, types and the date style on column D).
$ws.Rows(63).Copy()
$ws.Rows(65).Insert()

# 2) Duplicate the existing row 64 down to become the new row 66.
$ws.Rows(64).Copy()
$ws.Rows(66).Insert()

# 3) Overwrite row 63 with the new "Primera" record for 2022-08-03.
$ws.Range("D63").Value = 44776
$ws.Range("J63").Value = 200
$ws.Range("K63").Value = 700
$ws.Range("L63").Value = 800
$ws.Range("M63").Value = 750
$ws.Range("P63").Value = 750

# 4) Overwrite row 64 with the new "Segunda" record for 2022-08-03.
$ws.Range("D64").Value = 44776
$ws.Range("E64").Value = 16
$ws.Range("F64").Value = 100112040
$ws.Range("G64").Value = "Cilantro"
$ws.Range("H64").Value = "Sin especificar"
$ws.Range("I64").Value = "Segunda"
$ws.Range("J64").Value = 150
$ws.Range("K64").Value = 600
$ws.Range("L64").Value = 600
$ws.Range("M64").Value = 600
$ws.Range("N64").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O64").Value = "Provincia de Diguillín"
$ws.Range("P64").Value = 600
$ws.Range("Q64").Value = 1
$ws.Range("R64").Value = "Hortaliza"
